# Trade #117 closed at 2026-02-17 09:28:15 - unknown UNKNOWN +0.000%
#
# Updates the Summary and Strategy Status roll-up figures for the
# MarketMaking strategy after the new trade closed, and appends the
# new trade record (row 118) to both the "All Trades" and
# "MarketMaking" trade logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Summary sheet roll-up numbers
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.64   # Current Capital
$summary.Range("B4").Value = 0.65      # Total P&L $
$summary.Range("B5").Value = 0.11      # Total P&L %
$summary.Range("B6").Value = 117       # Total Trades
$summary.Range("B8").Value = 44        # Losing Trades
$summary.Range("B9").Value = 44.44     # Win Rate %

# ---------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.64     # Capital
$status.Range("D4").Value = 117        # Trades
$status.Range("E4").Value = 0.65       # P&L $
$status.Range("F4").Value = 0.64       # P&L %
$status.Range("G4").Value = 44.44      # Win Rate %

# ---------------------------------------------------------------
# 3. Append new trade row (118) to "All Trades" and "MarketMaking"
# ---------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Cells.Item(118, 1).Value = 117
    # Date column: force text so "2026-02-17" isn't auto-converted to a
    # date serial number, matching the literal-string storage used by
    # every other row in this column.
    $ws.Cells.Item(118, 2).Value = "'2026-02-17"
    $ws.Cells.Item(118, 3).Value = "09:28:09"
    $ws.Cells.Item(118, 4).Value = "MarketMaking"
    $ws.Cells.Item(118, 5).Value = "DOWN"
    $ws.Cells.Item(118, 6).Value = 0.91
    $ws.Cells.Item(118, 7).Value = 0.88
    $ws.Cells.Item(118, 8).Value = "CLOSED"
    $ws.Cells.Item(118, 9).Value = -3.2967
    $ws.Cells.Item(118, 10).Value = -0.03
    $ws.Cells.Item(118, 11).Value = 100.64
    $ws.Cells.Item(118, 12).Value = 0
    $ws.Cells.Item(118, 13).Value = 0
    $ws.Cells.Item(118, 14).Value = 0.6
    $ws.Cells.Item(118, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(118, 16).Value = "early_exit"
    $ws.Cells.Item(118, 17).Value = 0.17
}
